# Menus, Updated Icons, Documentation
# - Simplify the food icons: fill in real ingredient / allergen / diet info
#   for the "Apple Turnover" and "Peanut Butter Brownie" rows (previously
#   all placeholder "needed" text), and drop the explicit per-cell style
#   that was only ever re-stating the default formatting.
# - Shrink the backing table / selection back down to the actual data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ("Apple Turnover"): real ingredient list for the Ingredients column ---
$ws.Range("B2").Value2 = "Apple Filling  / Enriched Flour / Butter / Whole Eggs / Cinnamon / Salt"

# --- Row 6 ("Peanut Butter Brownie"): ingredients / allergens / local ingredients / diet ---
$ws.Range("B6").Value2 = "Peanut Butter /  Enriched Flour /  Margarine / White Chocolate / Chocolate / Skim Milk Powder / Whole Eggs /  Peanuts / Sugar"
$ws.Range("C6").Value2 = "Wheat, milk, eggs, soy, peanuts. May contain treenuts, sulphites."
$ws.Range("D6").Value2 = "NA"
$ws.Range("E6").Value2 = "VEG"

# --- Strip the redundant explicit style off all the "needed" placeholder cells ---
# (columns B:F, rows 2-16) so they fall back to the workbook's default style.
$ws.Range("B2:F16").Style = "Normal"

# --- Shrink the table back down to the real data range ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G16"))

# --- Move the selection/cursor ---
$ws.Range("E21").Select()
